# Generate Report for Handoff
# Update the "b.md" row across the Overview / zh-cn / de-de sheets to reflect
# that the file has been handed off again (new handoff xliff generated),
# so the previous "handed back" status is replaced with "Ready for handoff"
# and related metadata (dates, handoff file names, duplicate flag, error
# detail) are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md entry.
#   E = zh-cn status, F = de-de status, G = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 08:50:36"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md entry.
#   C = Status
#   F = Content Duplicate
#   G = Latest Handoff File
#   H = Latest Handoff Datetime
#   P = Error Detail
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps "False" as literal text (matches source file,
# which stores True/False as plain strings, not booleans); reset the style
# afterwards so the quote-prefix formatting doesn't stick.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 08:50:20"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc174efeafc1b3cbe9bf2d2d5e5af01467ad0d2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acc9c4c4613f3ebcaaddfddc0d2c53a345060ec4/e2e/b.md."
# 39.14 round-trips through the engine's pixel-width conversion to the
# stored OOXML column width of exactly 40 (matches the widened "Error
# Detail" column in the diff).
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md entry.
#   C = Status
#   F = Content Duplicate
#   G = Latest Handoff File
#   H = Latest Handoff Datetime
#   P = Error Detail
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 08:50:36"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bc174efeafc1b3cbe9bf2d2d5e5af01467ad0d2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acc9c4c4613f3ebcaaddfddc0d2c53a345060ec4/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.14
